$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Toggle which "D" column (spreadsheet column D, 4th col) values are
#    present vs. missing for a handful of RM rows (re-rolled missing mask).
$ws.Range("D5").ClearContents()
$ws.Range("D8").Value = -13.9
$ws.Range("D12").ClearContents()
$ws.Range("D14").Value = -13.1
$ws.Range("D18").ClearContents()

# 2. Drop the "RM 232" record entirely (row 26).
$ws.Rows.Item(26).Delete()

# 3. Drop the "SC 92" record entirely - after the deletion above it now sits
#    at row 27 (it used to be row 28).
$ws.Rows.Item(27).Delete()

# 4. Re-roll the missing mask on a few of the remaining SC rows, now shifted
#    into their final positions.
$ws.Range("B26").Value = -20.2   # SC 5   : was missing, now populated
$ws.Range("B27").ClearContents() # SC 101 : was populated, now missing
$ws.Range("C33").Value = 10.4    # SC 232 : was missing, now populated
